$d = $word.ActiveDocument

# The release-notes date reads "September 12, 2025" and needs to become
# "November 2025". The original text is split across two runs
# ("September 12" and ", 2025"); the target keeps a two-run split too
# ("November" and " 2025"). Find the whole phrase, clear it, then insert
# the two replacement pieces back-to-back so the run split is preserved
# (a single Find/Replace merges adjacent same-formatted runs into one).
$rng = $d.Content
$rng.Find.Execute("September 12, 2025")
$rng.Text = ""
$rng.Collapse(1)            # wdCollapseStart
$rng.InsertAfter("November")
$rng.Collapse(0)            # wdCollapseEnd
$rng.InsertAfter(" 2025")
